# Requests.xlsx - "Resquests and Responses done"
#
# Adds a new "User" request row (GET user/userID) and normalizes a couple of
# cell alignments (B5 / B11 go from vertical-only centering to full
# horizontal+vertical centering, matching the alignment already used by the
# other section-header cells in column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row: User / GET / user/userID -------------------------------
# Row 16 already holds "User" in B16; fill in the Resource/Request columns
# that were previously empty for that row.
$ws.Range("C16").Value = "GET"
$ws.Range("D16").Value = "user/userID"

# --- Alignment fix-ups ------------------------------------------------------
# B5 ("News") and B11 ("Product") were vertical-center only; align them with
# the rest of the section headers (horizontal + vertical center).
$xlCenter = -4108

$fixAlign = @("B5", "B11")
foreach ($addr in $fixAlign) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = $xlCenter
    $rng.VerticalAlignment = $xlCenter
}

# --- Last active selection --------------------------------------------------
$ws.Range("H7").Select()
